# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the notes master only)
#   ppt/theme/theme2.xml -> "Integral"     (used by the slide master / the
#                                            presentation's active design)
# The target revision swaps their contents, so the design that is actually
# applied to the slides/slide master switches from the "Integral" palette to
# the stock "Office" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# The object model only exposes the *active* theme (the one hanging off the
# slide master / Design) for editing, via ThemeColorScheme.Colors(i).RGB, so
# we push the "Office" theme's twelve colours onto it, in the fixed
# dk1,lt1,dk2,lt2,accent1..accent6,hlink,folHlink order.

function ToBGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$tcs = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = ToBGR($officeThemeColors[$i - 1])
}

# Best-effort: keep the theme/colour-scheme display names in step with the
# palette that is now applied (harmless if the host treats these as
# read-only).
try { $theme.Name = "Office Theme" } catch { }
try { $tcs.Name = "Office" } catch { }
